# Update the 2-Experts average scores in Sheet1 (res/SEED) to the latest run.
# Commit: "add RELU function. 2021/03/31 19:02"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated scores (column B) ---
$ws.Range("B24").Value = 0.9314
$ws.Range("B30").Value = 0.7652
$ws.Range("B31").Value = 0.953
$ws.Range("B32").Value = 0.7876
$ws.Range("B33").Value = 0.9025

# Row 40 flips from a "Bad" score to a "Good" one, so re-apply the
# conditional-style cell style along with the new value.
$ws.Range("B40").Value = 0.9509
$ws.Range("B40").Style = "好"

$ws.Range("B41").Value = 0.8837
$ws.Range("B43").Value = 0.7782

# B47 holds =AVERAGE(B2:B46); let it recalc against the new inputs.
$wb.Application.Calculate()
